$wb = $excel.ActiveWorkbook

$wsMeans = $wb.Worksheets.Item("Means")
$wsSD    = $wb.Worksheets.Item("Standard Deviations")

# --- Update headers: add "Within 5 miles..." and "Within 10 miles..." columns ---
$wsMeans.Range("F1").Value = "Within 5 miles of HFC production facility"
$wsMeans.Range("G1").Value = "Within 10 miles of HFC production facility"

$wsSD.Range("F1").Value = "Within 5 mile of HFC production facility SD"
$wsSD.Range("G1").Value = "Within 10 mile of HFC production facility SD"

# --- Means sheet: new F/G column values ---
$meansF = @(81, 12, 7.2, 3.8, 56, 9.9, 8.2, 20, 0.25)
$meansG = @(81, 9.1, 9.7, 5.8, 59, 7.7, 7.4, 20, 0.24)

for ($i = 0; $i -lt $meansF.Length; $i++) {
    $row = $i + 2
    $wsMeans.Cells.Item($row, 6).Value = $meansF[$i]
    $wsMeans.Cells.Item($row, 7).Value = $meansG[$i]
}

# --- Means sheet: corrected existing values (rows 9 and 10) ---
$wsMeans.Range("B9").Value = 26
$wsMeans.Range("D9").Value = 20
$wsMeans.Range("E9").Value = 21

$wsMeans.Range("B10").Value = 0.32
$wsMeans.Range("C10").Value = 0.22
$wsMeans.Range("D10").Value = 0.2
$wsMeans.Range("E10").Value = 0.25

# --- Standard Deviations sheet: new F/G column values ---
$sdF = @(24, 22, 9.3, 7.1, 25, 13, 11, 2.2, 0.05)
$sdG = @(22, 16, 14, 9.8, 25, 11, 11, 1.2, 0.05)

for ($i = 0; $i -lt $sdF.Length; $i++) {
    $row = $i + 2
    $wsSD.Cells.Item($row, 6).Value = $sdF[$i]
    $wsSD.Cells.Item($row, 7).Value = $sdG[$i]
}

# --- Standard Deviations sheet: corrected existing values (rows 9 and 10) ---
$wsSD.Range("B9").Value = 8.6
$wsSD.Range("D9").Value = 0
$wsSD.Range("E9").Value = 3.5

$wsSD.Range("C10").Value = 0.038
$wsSD.Range("D10").Value = 0
$wsSD.Range("E10").Value = 0.048
